# Update "Förändrad" date column (C) for rows 2-18 from 45205 to 45206 (2023-10-06 -> 2023-10-07)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
